$wb = $excel.ActiveWorkbook

# --- Sheet: tournament format ---
$ws1 = $wb.Worksheets.Item("tournament format")

# Row 5 (QR codes row): reword the directions text
$ws1.Range("C5").Value = 'Enter "Y" for QR codes to be included or "N" to be excluded.'
$ws1.Range("C5").WrapText = $true

# Row 6: new "text input" toggle row (was "Text"/"N"/old directions -> "text"/"Y"/new directions)
$ws1.Range("A6").Value = "text"
$ws1.Range("B6").Value = "Y"
$ws1.Range("C6").Value = 'Enter "Y" for text input to be included or "N" to be excluded.'
$ws1.Range("C6").WrapText = $true

# --- Sheet: text input ---
$ws6 = $wb.Worksheets.Item("text input")

$ws6.Range("B2").Value = "This is the placeholder text for the prelim schedule individualized for each team. Useful pieces of information to include in this section are: expected start time and length of time for lunch, information on tiebreakers, where to report back after lunch..."
$ws6.Range("B3").Value = "This is the placeholder text for the prelim schedule individualized for each room. Useful pieces of information to include in this section are: how to report protests, tiebreaker procedure, where to go for lunch…"

$ws6.Rows.Item(1).RowHeight = 40
$ws6.Rows.Item(2).RowHeight = 160
$ws6.Rows.Item(3).RowHeight = 160

$ws6.Range("B2").WrapText = $true
$ws6.Range("B3").WrapText = $true
$ws6.Range("B2:B3").HorizontalAlignment = -4108
$ws6.Range("B2:B3").VerticalAlignment = -4108
